# dynamic sink term added
#
# Sheet1 ("Cases") gains two new input columns ("T_whole" and "Xbed"),
# inserted right before the old "P_gas" column. The old "P_heat" /
# "nboil" pair of columns is collapsed into a single "P_fuel" column
# (keeping the P_heat value). H_bed and H_gap base-case values are also
# updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert two new columns at Z:AA (pushes P_gas/P_heat/nboil/... right by 2) ---
$ws.Range("Z1:AA1").EntireColumn.Insert()

# New column Z = T_whole
$ws.Cells.Item(1, 26).Value2 = "T_whole"
$ws.Cells.Item(2, 26).Value2 = "deg C"
$ws.Cells.Item(3, 26).Value2 = 500

# New column AA = Xbed
$ws.Cells.Item(1, 27).Value2 = "Xbed"
$ws.Cells.Item(2, 27).Value2 = "--"
$ws.Cells.Item(3, 27).Value2 = 0.6

# Column AB (28) is still the original "P_gas" column - untouched.

# Column AC (29) was "P_heat" - rename it to "P_fuel" (value/unit unchanged).
$ws.Cells.Item(1, 29).Value2 = "P_fuel"

# Column AD (30) was "nboil" - remove it entirely (merged into P_fuel).
$ws.Columns.Item(30).Delete()

# --- Update base-case values that changed ---
# H_bed base case: 0.65 -> 0.55
$ws.Cells.Item(3, 35).Value2 = 0.55
# H_gap base case: 0.5 -> 0.3
$ws.Cells.Item(3, 38).Value2 = 0.3

# --- Restore the view selection on Sheet1 ---
$ws.Range("AD1").Select()
